# Update the multiplication answer table to the newly generated values.
$d = $word.ActiveDocument

$replacements = @(
    @{ old = "64×98=6272"; new = "34×93=3162" },
    @{ old = "38×67=2546"; new = "71×57=4047" },
    @{ old = "41×16=656";  new = "12×92=1104" },
    @{ old = "44×58=2552"; new = "66×64=4224" },
    @{ old = "63×34=2142"; new = "73×57=4161" },
    @{ old = "36×42=1512"; new = "24×39=936" },
    @{ old = "62×52=3224"; new = "46×18=828" },
    @{ old = "33×11=363";  new = "29×86=2494" },
    @{ old = "40×66=2640"; new = "72×71=5112" },
    @{ old = "13×54=702";  new = "18×70=1260" },
    @{ old = "47×41=1927"; new = "58×80=4640" },
    @{ old = "83×16=1328"; new = "45×23=1035" },
    @{ old = "31×19=589";  new = "15×30=450" },
    @{ old = "60×89=5340"; new = "41×66=2706" },
    @{ old = "13×66=858";  new = "13×73=949" },
    @{ old = "21×77=1617"; new = "98×41=4018" },
    @{ old = "33×95=3135"; new = "23×73=1679" },
    @{ old = "69×63=4347"; new = "87×78=6786" },
    @{ old = "12×33=396";  new = "64×88=5632" },
    @{ old = "96×79=7584"; new = "97×48=4656" },
    @{ old = "53×82=4346"; new = "68×43=2924" },
    @{ old = "46×82=3772"; new = "75×19=1425" },
    @{ old = "23×85=1955"; new = "78×72=5616" },
    @{ old = "97×45=4365"; new = "78×53=4134" },
    @{ old = "71×33=2343"; new = "56×60=3360" }
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $pair.new, 2)
}
